$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D1 header holds the shared string "O(n.m)"; change it to "O(m)" ---
$ws.Range("D1").Value2 = "O(m)"

# --- Column D: complexity was A*B (n.m); it is now just B (m) ---
$ws.Range("D2").Formula = "=B2"
$ws.Range("D3:D31").Formula = "=B3"

# --- Refresh the two "Complejidad O(...)" charts so their title / series name show O(m) ---
$charts = $ws.ChartObjects()
for ($i = 1; $i -le $charts.Count; $i++) {
    $co = $charts.Item($i)
    $chart = $co.Chart
    $ser = $chart.SeriesCollection().Item(1)
    if ($ser.Name -like "O(*") {
        if ($chart.HasTitle) {
            $title = $chart.ChartTitle.Text
            $chart.ChartTitle.Text = $title -replace "O\(n.m\)", "O(m)"
        }
        $ser.Name = "O(m)"
    }
}

# --- Sheet view: zoom to 85%, drop the old frozen scroll position, select V11 ---
$win = $excel.ActiveWindow
$win.Zoom = 85
[void]$ws.Range("V11").Select()
